$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B=question, C=answer shift right)
$ws.Columns("B:B").Insert()

# New header + values for the inserted "relevant_chunks_id" column
$ws.Range("B1").Value = "relevant_chunks_id"
$ws.Range("B2").Value = "id1,id2,id3"
$ws.Range("B3").Value = "id1,id2,id4"
$ws.Range("B4").Value = "id1,id2,id5"
$ws.Range("B5").Value = "id1,id2,id6"
$ws.Range("B6").Value = "id1,id2,id7"
$ws.Range("B7").Value = "id1,id2,id8"
$ws.Range("B8").Value = "id1,id2,id9"
$ws.Range("B9").Value = "id1,id2,id10"
$ws.Range("B10").Value = "id1,id2,id11"
$ws.Range("B11").Value = "id1,id2,id12"
$ws.Range("B12").Value = "id1,id2,id13"

# Column width for the new column
$ws.Columns("B:B").ColumnWidth = 12.5

# Style the new data cells (not bold, centered horizontally/vertically)
$ws.Range("B2:B12").Font.Bold = $false
$ws.Range("B2:B12").HorizontalAlignment = -4108
$ws.Range("B2:B12").VerticalAlignment = -4108
